# Regenerate save_data to use K (strike count) instead of Strike# (old raw value).
# The recomputed K values (column G) replace the previously stored values for
# every data row (rows 2-40) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 5
    4  = 4
    5  = 5
    6  = 10
    7  = 4
    8  = 4
    9  = 4
    10 = 7
    11 = 5
    12 = 5
    13 = 2
    14 = 6
    15 = 9
    16 = 10
    17 = 11
    18 = 8
    19 = 5
    20 = 9
    21 = 10
    22 = 5
    23 = 8
    24 = 6
    25 = 5
    26 = 9
    27 = 5
    28 = 8
    29 = 12
    30 = 8
    31 = 5
    32 = 8
    33 = 7
    34 = 7
    35 = 7
    36 = 6
    37 = 4
    38 = 8
    39 = 5
    40 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
